# Revert "Merge pull request #48 from LakeFishing/main"
# - Cell A3 on "Sheet" changes from "投籃" back to "投藍"
# - Selection moves from A3 back to G5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$ws.Range("A3").Value = "投藍"

$ws.Range("G5").Select()
